$d = $word.ActiveDocument

# The Title, Author and Date paragraphs at the top of the document were
# originally split into a separate <w:r> run for every word and every
# in-between space. Re-flatten each of those paragraphs down to a single
# run that carries the whole sentence (what Word ends up with once the
# paragraph's text is rewritten as a whole).

function Merge-ParagraphRuns($paragraph) {
    $range = $paragraph.Range
    # Exclude the trailing paragraph mark from the range we rewrite.
    $range.MoveEnd(1, -1) | Out-Null
    $finalText = $range.Text

    # Writing the same text straight back is a no-op, so stash a
    # placeholder first to force a real content change - that collapses
    # the run-per-word fragments into a single run - then write the real
    # text into that single run.
    $range.Text = "@@PLACEHOLDER@@"

    $range2 = $paragraph.Range
    $range2.MoveEnd(1, -1) | Out-Null
    $range2.Text = $finalText
}

$targetStyles = @("Title", "Author", "Date")
foreach ($p in $d.Paragraphs) {
    if ($targetStyles -contains $p.Style.NameLocal) {
        Merge-ParagraphRuns $p
    }
}

Write-Output "Title : $($d.Paragraphs.Item(1).Range.Text)"
Write-Output "Author: $($d.Paragraphs.Item(2).Range.Text)"
Write-Output "Date  : $($d.Paragraphs.Item(3).Range.Text)"
